$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right count and Wrong penalty corrections
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": totals recomputed for the marks error fix
$ws.Range("B12").Value = 104
$ws.Range("C12").Value = -4
$ws.Range("E12").Value = "100 / 112"
